$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "DBServer_1"
$ws.Range("C16").Value = "DBServer"
$ws.Range("C16").Select()
